# Weekly update: insert the latest price observations (date serial 45146)
# for "Betarraga" at "Vega Monumental Concepción" into the historical
# series. The new week's rows are inserted right above the most recent
# existing observation (row 406), shifting the rest of the table down by
# two rows (the table keeps growing, so the two oldest/ tail rows get
# duplicated at the bottom automatically by the row insert - this mirrors
# how the source data is appended weekly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 406/407; everything from 406 downward
# shifts down by two rows (406->408, 407->409, ..., 545->547).
$ws.Rows("406:407").Insert()

# New row 406: "Primera" quality for the new week (date 45146)
$ws.Range("A406").Value = 11
$ws.Range("B406").Value = "Vega Monumental Concepción"
$ws.Range("C406").Value = "Bíobío"
$ws.Range("D406").Value = 45146
$ws.Range("E406").Value = 8
$ws.Range("F406").Value = 100114014
$ws.Range("G406").Value = "Betarraga"
$ws.Range("H406").Value = "Sin especificar"
$ws.Range("I406").Value = "Primera"
$ws.Range("J406").Value = 200
$ws.Range("K406").Value = 600
$ws.Range("L406").Value = 700
$ws.Range("M406").Value = 650
$ws.Range("N406").Value = "`$/paquete 5 unidades"
$ws.Range("O406").Value = "Región Metropolitana"
$ws.Range("P406").Value = 130
$ws.Range("Q406").Value = 5
$ws.Range("R406").Value = "Hortaliza"

# New row 407: "Segunda" quality for the new week (date 45146)
$ws.Range("A407").Value = 11
$ws.Range("B407").Value = "Vega Monumental Concepción"
$ws.Range("C407").Value = "Bíobío"
$ws.Range("D407").Value = 45146
$ws.Range("E407").Value = 8
$ws.Range("F407").Value = 100114014
$ws.Range("G407").Value = "Betarraga"
$ws.Range("H407").Value = "Sin especificar"
$ws.Range("I407").Value = "Segunda"
$ws.Range("J407").Value = 100
$ws.Range("K407").Value = 500
$ws.Range("L407").Value = 500
$ws.Range("M407").Value = 500
$ws.Range("N407").Value = "`$/paquete 5 unidades"
$ws.Range("O407").Value = "Región Metropolitana"
$ws.Range("P407").Value = 100
$ws.Range("Q407").Value = 5
$ws.Range("R407").Value = "Hortaliza"
